{"js": "// Split three long \"run-on\" paragraphs (Programa PT, Programa EN, Bibliografia)\n// into multiple <w:t> segments separated by <w:br/> line breaks, matching the\n// topic/reference boundaries the author introduced in the target revision.\n// Implemented generically (search the paragraph text for marker substrings)\n// rather than by hard-coded paragraph index, then rebuilt with insertOoxml so\n// the break stays inside the SAME run as the surrounding text (<w:r><w:t/>\n// <w:br/><w:t/>...</w:r>), exactly like the diff.\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\nfunction splitOnMarkers(text, markers) {\n  const idxs = [];\n  let pos = 0;\n  for (const mk of markers) {\n    const i = text.indexOf(mk, pos);\n    if (i === -1) {\n      throw new Error(\"marker not found: \" + mk);\n    }\n    idxs.push(i);\n    pos = i + 1;\n  }\n  const segs = [];\n  let prev = 0;\n  for (const i of idxs) {\n    segs.push(text.slice(prev, i));\n    prev = i;\n  }\n  segs.push(text.slice(prev));\n  return segs;\n}\n\nfunction buildParagraphOoxml(segs, italic) {\n  const rPr = italic ? \"<w:rPr><w:i/></w:rPr>\" : \"\";\n  const parts = [];\n  for (let i = 0; i < segs.length; i++) {\n    const preserve = /^\\s|\\s$/.test(segs[i]) ? ' xml:space=\"preserve\"' : \"\";\n    parts.push(`<w:t${preserve}>${escapeXml(segs[i])}</w:t>`);\n    if (i !== segs.length - 1) {\n      parts.push(\"<w:br/>\");\n    }\n  }\n  const run = `<w:r>${rPr}${parts.join(\"\")}</w:r>`;\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    \"</Relationships>\" +\n    \"</pkg:xmlData></pkg:part>\" +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    `<w:body><w:p>${run}</w:p></w:body></w:document>` +\n    \"</pkg:xmlData></pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Each entry: a substring that uniquely identifies the target paragraph's\n// start, the ordered markers that begin each new segment (first marker is\n// skipped since it is the start of the text itself).\nconst targets = [\n  {\n    startsWith: \"Integra\u00e7\u00e3o de fun\u00e7\u00f5es reais: Primitivas\",\n    markers: [\n      \"O espa\u00e7o euclidiano R^n:\",\n      \"Fun\u00e7\u00f5es de n v\u00e1rias vari\u00e1veis Reais:\",\n      \"Limites e Continuidade:\",\n      \"Diferenciabilidade:\",\n      \"M\u00e1ximos e m\u00ednimos:\",\n    ],\n  },\n  {\n    startsWith: \"Integration of real functions: Primitive function\",\n    markers: [\n      \"The Euclidian Espace R^n:\",\n      \"Function of n Real variables:\",\n      \"Limits and continuity:\",\n      \"Differentiability:\",\n      \"Maximum and Minimum:\",\n    ],\n  },\n  {\n    startsWith: \"GUIDORIZZI, Hamilton L.\",\n    markers: [\n      \"LEITHOLD, Louis.\",\n      \"ANTON, Howard\",\n      \"SIMMONS, George F.\",\n      \"STEWART, James.\",\n      \"THOMAS, George B.\",\n    ],\n  },\n];\n\nfor (const item of paragraphs.items) {\n  item.load(\"text\");\n}\nawait context.sync();\n\nfor (const target of targets) {\n  const para = paragraphs.items.find((p) => p.text.indexOf(target.startsWith) === 0);\n  if (!para) {\n    throw new Error(\"paragraph not found for: \" + target.startsWith);\n  }\n  para.font.load(\"italic\");\n  await context.sync();\n\n  const segments = splitOnMarkers(para.text, target.markers);\n  const ooxml = buildParagraphOoxml(segments, para.font.italic === true);\n  para.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Split three long \"run-on\" paragraphs (Programa PT, Programa EN, Bibliografia)\n# into multiple <w:t> segments separated by <w:br/> line breaks, matching the\n# topic/reference boundaries the author introduced in the target revision.\n#\n# Paragraphs are located generically by a leading substring (not a hard-coded\n# index). Each target paragraph's Range is rebuilt with Range.InsertXML(...)\n# using an explicit flat-OPC WordOpenXML fragment, so the break lands inside\n# the SAME run as the surrounding text (<w:r><w:t/><w:br/><w:t/>...</w:r>),\n# with explicit xml:space=\"preserve\" exactly where the original text has\n# leading/trailing whitespace on a segment -- exactly like the diff -- and\n# the run keeps its original formatting (e.g. italics) because we read it\n# from Range.Font.Italic before rebuilding.\n\nfunction Split-OnMarkers($text, $markers) {\n  $idxs = @()\n  $pos = 0\n  foreach ($mk in $markers) {\n    $i = $text.IndexOf($mk, $pos)\n    if ($i -eq -1) {\n      throw \"marker not found: $mk\"\n    }\n    $idxs += $i\n    $pos = $i + 1\n  }\n  $segs = @()\n  $prev = 0\n  foreach ($i in $idxs) {\n    $segs += $text.Substring($prev, $i - $prev)\n    $prev = $i\n  }\n  $segs += $text.Substring($prev)\n  return $segs\n}\n\nfunction Escape-Xml($s) {\n  $s = $s.Replace(\"&\", \"&amp;\")\n  $s = $s.Replace(\"<\", \"&lt;\")\n  $s = $s.Replace(\">\", \"&gt;\")\n  $s = $s.Replace('\"', \"&quot;\")\n  return $s\n}\n\nfunction Build-ParagraphOoxml($segments, $italic) {\n  $rPr = \"\"\n  if ($italic) {\n    $rPr = \"<w:rPr><w:i/></w:rPr>\"\n  }\n  $parts = @()\n  for ($i = 0; $i -lt $segments.Length; $i++) {\n    $seg = $segments[$i]\n    $preserve = \"\"\n    if ($seg.Length -gt 0 -and (($seg.Substring(0,1) -match '\\s') -or ($seg.Substring($seg.Length-1,1) -match '\\s'))) {\n      $preserve = ' xml:space=\"preserve\"'\n    }\n    $parts += (\"<w:t\" + $preserve + \">\" + (Escape-Xml $seg) + \"</w:t>\")\n    if ($i -ne $segments.Length - 1) {\n      $parts += \"<w:br/>\"\n    }\n  }\n  $run = \"<w:r>\" + $rPr + ([string]::Join(\"\", $parts)) + \"</w:r>\"\n  $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    ('<w:body><w:p>' + $run + '</w:p></w:body></w:document>') +\n    '</pkg:xmlData></pkg:part>' +\n    '</pkg:package>'\n  return $xml\n}\n\n$d = $word.ActiveDocument\n\n$targets = @(\n  @{\n    StartsWith = \"Integra\u00e7\u00e3o de fun\u00e7\u00f5es reais: Primitivas\"\n    Markers = @(\n      \"O espa\u00e7o euclidiano R^n:\",\n      \"Fun\u00e7\u00f5es de n v\u00e1rias vari\u00e1veis Reais:\",\n      \"Limites e Continuidade:\",\n      \"Diferenciabilidade:\",\n      \"M\u00e1ximos e m\u00ednimos:\"\n    )\n  },\n  @{\n    StartsWith = \"Integration of real functions: Primitive function\"\n    Markers = @(\n      \"The Euclidian Espace R^n:\",\n      \"Function of n Real variables:\",\n      \"Limits and continuity:\",\n      \"Differentiability:\",\n      \"Maximum and Minimum:\"\n    )\n  },\n  @{\n    StartsWith = \"GUIDORIZZI, Hamilton L.\"\n    Markers = @(\n      \"LEITHOLD, Louis.\",\n      \"ANTON, Howard\",\n      \"SIMMONS, George F.\",\n      \"STEWART, James.\",\n      \"THOMAS, George B.\"\n    )\n  }\n)\n\nforeach ($target in $targets) {\n  $para = $null\n  foreach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith($target.StartsWith)) {\n      $para = $p\n      break\n    }\n  }\n  if ($para -eq $null) {\n    throw \"paragraph not found for: \" + $target.StartsWith\n  }\n\n  $italic = [bool]$para.Range.Font.Italic\n\n  # Original paragraph text includes the trailing paragraph mark; strip it\n  # before splitting so it is not glued onto the last segment.\n  $full = $para.Range.Text\n  if ($full.EndsWith([char]13)) {\n    $full = $full.Substring(0, $full.Length - 1)\n  }\n\n  $segments = Split-OnMarkers $full $target.Markers\n  $ooxml = Build-ParagraphOoxml $segments $italic\n\n  [void]$para.Range.InsertXML($ooxml)\n}\n"}
